$d = $word.ActiveDocument
$endPos = $d.Content.End
$r = $d.Range($endPos, $endPos)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="es-AR"/></w:rPr><w:t>¿Qué cambia que el enunciado diga que cuando tengo que borrar un elemento de la lista puede que no exista o seguro existe?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">Si el elemento puede no existir: hay que preguntar por </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>act</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> &lt;&gt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>nil</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">, tanto en el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>while</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> que recorre como en el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>if</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> que elimina</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">Si el elemento seguro existe: no hace falta preguntar por </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>act</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> &lt;&gt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>nil</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> ni en el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>while</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> ni en el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>if</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
